# Updates cryptos list values (price + 1h volume change) and reorders a few rows
# to match the refreshed data pulled on Sun Dec 10 04:29:36 UTC 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.883.31"
$ws.Range("E2").Value = "  -0.97%  "

# Row 3
$ws.Range("D3").Value = "2.351.74"
$ws.Range("E3").Value = "  -0.51%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.674"
$ws.Range("E5").Value = "  -2.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.59"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.24"
$ws.Range("E7").Value = "  -1.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +0.48%  "

# Row 10
$ws.Range("E10").Value = "  -2.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.08"
$ws.Range("E11").Value = "  +2.16%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.52"
$ws.Range("E12").Value = "  +5.85%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.36"
$ws.Range("E13").Value = "  -1.87%  "

# Row 14
$ws.Range("E14").Value = "  +0.33%  "

# Row 15
$ws.Range("D15").Value = "2.701.97"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.43"
$ws.Range("E16").Value = "  -3.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.909"
$ws.Range("E17").Value = "  -0.80%  "

# Row 18
$ws.Range("D18").Value = "2.350.40"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("D19").Value = "43.789.63"
$ws.Range("E19").Value = "  -1.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("E20").Value = "  -1.81%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.70"
$ws.Range("E21").Value = "  -0.70%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.45"
$ws.Range("E22").Value = "  -1.40%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.99"
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  +17.96%  "

# Row 25
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.73"
$ws.Range("E26").Value = "  -0.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("E27").Value = "  -3.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.57"
$ws.Range("E28").Value = "  -2.36%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.67"
$ws.Range("E30").Value = "  -0.15%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.80"
$ws.Range("E31").Value = "  +1.46%  "

# Row 32
$ws.Range("E32").Value = "  -0.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  +0.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").Value = "  -0.20%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  -3.85%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.43"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  -1.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"
$ws.Range("E38").Value = "  -3.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.37"
$ws.Range("E39").Value = "  -2.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0279"
$ws.Range("E40").Value = "  +1.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "66.81"
$ws.Range("E41").Value = "  +25.36%  "

# Row 42
$ws.Range("E42").Value = "  +11.37%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.07"
$ws.Range("E43").Value = "  +13.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.12"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.10"
$ws.Range("E45").Value = "  -0.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").Value = "  +1.97%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.26"
$ws.Range("E47").Value = "  -0.20%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.04%  "

# Row 50
$ws.Range("E50").Value = "  -1.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.16"
$ws.Range("E51").Value = "  -2.45%  "
